$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.651.03"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.069.30"
$ws.Range("E3").Value = "  -3.73%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.69"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.52"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.066.47"
$ws.Range("E9").Value = "  -3.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  -4.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.85"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -3.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.71"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.574.39"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.584.67"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.10"
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.066.04"
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.03"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.24"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.46"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.29"
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.75"
$ws.Range("E26").Value = "  -4.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.43"
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.65"
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.14"
$ws.Range("E32").Value = "  -5.67%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").Value = "  -8.34%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.96"
$ws.Range("E34").Value = "  -5.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0817"
$ws.Range("E35").Value = "  -5.55%  "
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.96"
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  -4.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  -5.05%  "
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.14"
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "436.46"
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.286"
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.60"
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0358"
$ws.Range("E46").Value = "  -5.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.796.25"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.55"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.08"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.21"
$ws.Range("E51").Value = "  -3.10%  "
